$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 4954
$ws.Range("F4").Value = 1871
$ws.Range("F6").Value = 38
$ws.Range("F7").Value = 87
$ws.Range("F8").Value = 91
$ws.Range("F9").Value = 242
$ws.Range("F11").Value = 1100
$ws.Range("F12").Value = 367
$ws.Range("F14").Value = 61
$ws.Range("F15").Value = 110
$ws.Range("F16").Value = 17
$ws.Range("F17").Value = 236
$ws.Range("F18").Value = 132
$ws.Range("F19").Value = 84
$ws.Range("F20").Value = 1218
$ws.Range("F21").Value = 472
$ws.Range("F22").Value = 157
$ws.Range("F23").Value = 267
$ws.Range("F25").Value = 575
$ws.Range("F26").Value = 1027
$ws.Range("F27").Value = 54
$ws.Range("F28").Value = 1930
$ws.Range("F29").Value = 2387
$ws.Range("F30").Value = 1165
$ws.Range("F32").Value = 99
$ws.Range("F33").Value = 323
$ws.Range("F34").Value = 399
$ws.Range("F35").Value = 719
$ws.Range("F36").Value = 690
$ws.Range("F37").Value = 96
$ws.Range("F39").Value = 716
$ws.Range("F40").Value = 142
$ws.Range("F41").Value = 528
$ws.Range("F42").Value = 600
$ws.Range("F43").Value = 273
$ws.Range("F44").Value = 186

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value = 27
$ws.Range("F15").Value = 237
$ws.Range("F22").Value = 11
$ws.Range("F24").Value = 34

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 843

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 843
$ws.Range("F3").Value = 4954
$ws.Range("F4").Value = 1871
$ws.Range("F5").Value = 38
$ws.Range("F7").Value = 87
$ws.Range("F9").Value = 27
$ws.Range("F10").Value = 91
$ws.Range("F11").Value = 242
$ws.Range("F14").Value = 1100
$ws.Range("F15").Value = 367
$ws.Range("F17").Value = 61
$ws.Range("F18").Value = 110
$ws.Range("F19").Value = 236
$ws.Range("F21").Value = 132
$ws.Range("F22").Value = 85
$ws.Range("F23").Value = 1218
$ws.Range("F24").Value = 472
$ws.Range("F25").Value = 157
$ws.Range("F26").Value = 267
$ws.Range("F28").Value = 1027
$ws.Range("F29").Value = 1930
$ws.Range("F30").Value = 2387
$ws.Range("F32").Value = 1165
$ws.Range("F36").Value = 99
$ws.Range("F37").Value = 323
$ws.Range("F38").Value = 399
$ws.Range("F41").Value = 719
$ws.Range("F42").Value = 690
$ws.Range("F43").Value = 716
$ws.Range("F44").Value = 142
$ws.Range("F45").Value = 528
$ws.Range("F46").Value = 600
$ws.Range("F47").Value = 273
$ws.Range("F48").Value = 187
$ws.Range("F49").Value = 34
